$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Previous statement periods (rows 16-20, column E) are replaced with the
# new/updated list of periods, in reverse order (most recent period first),
# mirroring "Elimina EC anteriores y se agregan nuevos" (old balances removed,
# new ones added). Column F (Valor Mora) follows the same period -> amount
# mapping as before: period 2109 -> 40000, all other periods -> 48000.
$ws.Range("E16").Value = "2201"
$ws.Range("E17").Value = "2112"
$ws.Range("E18").Value = "2111"
$ws.Range("E19").Value = "2110"
$ws.Range("E20").Value = "2109"

$ws.Range("F16").Value = 40000
$ws.Range("F17").Value = 48000
$ws.Range("F18").Value = 48000
$ws.Range("F19").Value = 48000
$ws.Range("F20").Value = 48000
